$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 165 - pushes old rows 165..190 down to 166..191.
$ws.Rows.Item(165).Insert()

# Restore per-cell number/cell formatting (the blank Insert() drops borders /
# normalises styles) by pulling the exact formats already used elsewhere in
# this same column from the original sheet.
$ws.Range("A7").Copy()
$ws.Range("A165").PasteSpecial(-4122)

$ws.Range("B16").Copy()
$ws.Range("B165").PasteSpecial(-4122)

$ws.Range("C30").Copy()
$ws.Range("C165").PasteSpecial(-4122)

$ws.Range("D7").Copy()
$ws.Range("D165").PasteSpecial(-4122)

$ws.Range("E7").Copy()
$ws.Range("E165").PasteSpecial(-4122)

$ws.Range("F7").Copy()
$ws.Range("F165").PasteSpecial(-4122)

$ws.Range("G7").Copy()
$ws.Range("G165").PasteSpecial(-4122)

$ws.Range("H12").Copy()
$ws.Range("H165").PasteSpecial(-4122)

$ws.Range("I7").Copy()
$ws.Range("I165").PasteSpecial(-4122)

$ws.Range("J20").Copy()
$ws.Range("J165").PasteSpecial(-4122)

$ws.Range("K3").Copy()
$ws.Range("K165").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0

# Row height matches the new "SB" row (32pt, same as the BATS row above it).
$ws.Rows.Item(165).RowHeight = 32

# Populate the new Starbase (SB) line.
$ws.Range("A165").Value = "SB"
$ws.Range("B165").Value = "R1.01"
$ws.Range("C165").Value = "36(12)PP◆/`n18(6)P◆"
$ws.Range("D165").Value = "F&E"
$ws.Range("E165").Value = 10
$ws.Range("F165").Value = "Y140"
$ws.Range("G165").Value = "SB(1)"
$ws.Range("H165").Value = "From BATS: 30+6 From STB: 20+6"
$ws.Range("I165").Value = "Upgrade Only"
$ws.Range("J165").Value = 0
$ws.Range("K165").Value = "Starbase. See EW Chart. See (433.41), (441.0), and (510.3). Module`nLimits 4FTM+2PFM."
